$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B8").Value = "Deaths Door"
$ws.Range("E8").Value = 1
$ws.Range("G12").Select()
